# The "<id>p034v_1</id>" marker in the first <div> block is split across
# three separate runs:
#   run1: "<id>"   (Courier New, color 7f6000, sz/szCs 18)
#   run2: "p034v_1" (color 000000)
#   run3: "</id>"  (Courier New, color 7f6000, sz/szCs 18)
# The edit re-downloads/normalizes this marker into a single run that
# keeps run1's formatting and carries the full "<id>p034v_1</id>" text.
#
# (The sibling "<id>fig_p034v_1</id>" marker elsewhere in the document
# must NOT be touched, so we locate the exact, unique "<id>p034v_1</id>"
# span first.)

$d = $word.ActiveDocument

$marker = "<id>p034v_1</id>"
$prefix = "<id>"

# Locate the unique span covering all three runs.
$span = $d.Content
$span.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$spanStart = $span.Start
$spanEnd = $span.End

# Delete everything after the leading "<id>" (i.e. runs 2 and 3: the
# "p034v_1" and "</id>" runs), leaving run1's own text node alone so its
# formatting/markup (incl. xml:space) is preserved verbatim.
$tail = $d.Range($spanStart + $prefix.Length, $spanEnd)
$tail.Delete()

# Re-append the removed text ("p034v_1</id>") onto run1's range; Word
# extends the existing run rather than minting a differently-formatted
# one since the insertion point sits at the end of run1's own text.
$run1 = $d.Range($spanStart, $spanStart + $prefix.Length)
$run1.InsertAfter($marker.Substring($prefix.Length))
